# "maj perk data use" - update the rarity values used in the perk data
# table to the lowercase form ("epic"/"rare"/"uncommon" instead of
# "Epic"/"Rare"/"Uncommon"), remove the now-unneeded "READ FIRST"
# instructions textbox, and leave the sheet's selection where the author
# last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Rarity") values: lower-case the rarity labels.
$ws.Range("C2").Value = "epic"
$ws.Range("C3").Value = "epic"
$ws.Range("C4").Value = "rare"
$ws.Range("C5").Value = "rare"
$ws.Range("C6").Value = "uncommon"
$ws.Range("C7").Value = "uncommon"

# Remove the "READ FIRST" instructions textbox shape from the sheet.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Restore the author's final selected cell.
$ws.Range("E8").Select()
